# Regen save_data: switch the "K" column (strikeouts, header in G1) from the
# old Strike# derived count to the recalculated K value for each start.
# Row 2 (A=0) and row 49 (A=47) are unchanged - their recalculated K equals
# the value already on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    3  = 1
    4  = 1
    5  = 3
    6  = 0
    7  = 1
    8  = 3
    9  = 1
    10 = 0
    11 = 1
    12 = 4
    13 = 3
    14 = 0
    15 = 0
    16 = 1
    17 = 2
    18 = 0
    19 = 0
    20 = 1
    21 = 2
    22 = 2
    23 = 2
    24 = 1
    25 = 1
    26 = 0
    27 = 0
    28 = 1
    29 = 1
    30 = 2
    31 = 1
    32 = 2
    33 = 0
    34 = 1
    35 = 0
    36 = 0
    37 = 4
    38 = 1
    39 = 1
    40 = 1
    41 = 1
    42 = 0
    43 = 0
    44 = 2
    45 = 0
    46 = 1
    47 = 1
    48 = 1
    50 = 1
    51 = 2
    52 = 0
    53 = 2
    54 = 1
    55 = 1
    56 = 2
    57 = 1
    58 = 0
    59 = 1
    60 = 1
    61 = 0
    62 = 2
    63 = 1
    64 = 1
    65 = 2
    66 = 2
    67 = 0
    68 = 2
    69 = 1
    70 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}

Write-Output "Updated $($kValues.Count) K values in column G"
